# Normalize "Recorded By" (column G) entries: when the value starts with
# a leading "System," token, move that leading token to the end of the
# comma-separated list (swap the first and last segments), preserving the
# original casing of each token exactly as it appears.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value2

    if ($text -ne $null -and $text.StartsWith("System,")) {
        $parts = $text -split ", "
        $first = $parts[0]
        $lastIdx = $parts.Count - 1
        $last = $parts[$lastIdx]
        $parts[0] = $last
        $parts[$lastIdx] = $first
        $newText = $parts -join ", "
        $cell.Value2 = $newText
    }
}
